$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("shift")

# --- Defined names (workbook.xml <definedNames>) ---
# These are sheet-local names scoped to the "shift" sheet (index 5 -> localSheetId 4)
$ws.Names.Add("res", "=shift!`$D`$2:`$D`$16")
$ws.Names.Add("res_1", "=shift!`$E`$2:`$E`$16")
$ws.Names.Add("res_2", "=shift!`$F`$2:`$F`$16")

# --- Styles: apply a new font (black Calibri) across the used range, in the
# order that reproduces the target cellXfs layout (font-only, text-left-align,
# plain-left-align, numeric-format) ---
$ws.Range("B1:C1").Font.Color = 0
$ws.Range("E1:F1").Font.Color = 0
$ws.Range("C2:F16").Font.Color = 0
$ws.Range("A1:A16").Font.Color = 0
$ws.Range("D1").Font.Color = 0
$ws.Range("B2:B16").Font.Color = 0
$ws.Range("B2:B16").NumberFormat = "0.000000"

# --- Corrected data values for columns B-F, rows 2-16 ---
$ws.Range("B2").Value = 7.0226730000000002
$ws.Range("C2").Value = -1.1639764720000001
$ws.Range("D2").Value = -1.1201000000000001
$ws.Range("E2").Value = -1.1343000000000001
$ws.Range("F2").Value = -1.1306
$ws.Range("B3").Value = 6.9563560000000004
$ws.Range("C3").Value = -1.174777878
$ws.Range("D3").Value = -1.1435999999999999
$ws.Range("E3").Value = -1.1559999999999999
$ws.Range("F3").Value = -1.1534
$ws.Range("B4").Value = 6.8900379999999997
$ws.Range("C4").Value = -1.186285249
$ws.Range("D4").Value = -1.1645000000000001
$ws.Range("E4").Value = -1.1732
$ws.Range("F4").Value = -1.1712
$ws.Range("B5").Value = 6.8237209999999999
$ws.Range("C5").Value = -1.1951197140000001
$ws.Range("D5").Value = -1.1823999999999999
$ws.Range("E5").Value = -1.1875
$ws.Range("F5").Value = -1.1861999999999999
$ws.Range("B6").Value = 6.757403
$ws.Range("C6").Value = -1.201789497
$ws.Range("D6").Value = -1.1964999999999999
$ws.Range("E6").Value = -1.1984999999999999
$ws.Range("F6").Value = -1.1978
$ws.Range("B7").Value = 6.6910860000000003
$ws.Range("C7").Value = -1.205037084
$ws.Range("D7").Value = -1.2061999999999999
$ws.Range("E7").Value = -1.2055
$ws.Range("F7").Value = -1.2055
$ws.Range("B8").Value = 6.6247680000000004
$ws.Range("C8").Value = -1.204158794
$ws.Range("D8").Value = -1.2107000000000001
$ws.Range("E8").Value = -1.2079
$ws.Range("F8").Value = -1.2084999999999999
$ws.Range("B9").Value = 6.5584509999999998
$ws.Range("C9").Value = -1.1990502110000001
$ws.Range("D9").Value = -1.2091000000000001
$ws.Range("E9").Value = -1.2049000000000001
$ws.Range("F9").Value = -1.2060999999999999
$ws.Range("B10").Value = 6.4921329999999999
$ws.Range("C10").Value = -1.188868799
$ws.Range("D10").Value = -1.2003999999999999
$ws.Range("E10").Value = -1.1956
$ws.Range("F10").Value = -1.1974
$ws.Range("B11").Value = 6.4258160000000002
$ws.Range("C11").Value = -1.1728875329999999
$ws.Range("D11").Value = -1.1833
$ws.Range("E11").Value = -1.1791
$ws.Range("F11").Value = -1.1814
$ws.Range("B12").Value = 6.3594980000000003
$ws.Range("C12").Value = -1.1503113570000001
$ws.Range("D12").Value = -1.1567000000000001
$ws.Range("E12").Value = -1.1543000000000001
$ws.Range("F12").Value = -1.1569
$ws.Range("B13").Value = 6.2931809999999997
$ws.Range("C13").Value = -1.1203201810000001
$ws.Range("D13").Value = -1.119
$ws.Range("E13").Value = -1.1197999999999999
$ws.Range("F13").Value = -1.1226
$ws.Range("B14").Value = 6.2268629999999998
$ws.Range("C14").Value = -1.082144172
$ws.Range("D14").Value = -1.0686
$ws.Range("E14").Value = -1.0744
$ws.Range("F14").Value = -1.077
$ws.Range("B15").Value = 6.1605460000000001
$ws.Range("C15").Value = -1.0345742469999999
$ws.Range("D15").Value = -1.0038
$ws.Range("E15").Value = -1.0162
$ws.Range("F15").Value = -1.0184
$ws.Range("B16").Value = 6.0942280000000002
$ws.Range("C16").Value = -0.97661555099999997
$ws.Range("D16").Value = -0.92230000000000001
$ws.Range("E16").Value = -0.94359999999999999
$ws.Range("F16").Value = -0.94489999999999996

Write-Host "shift sheet corrected"
